$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so numeric-looking strings
# (e.g. "26.726.72", "1.380", "0.9973") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.726.72'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '1.731.28'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Value = '0.9973'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '0.9978'
$ws.Range("D7").Value = '0.4918'
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("D8").Value = '0.2619'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.06222'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = '1.729.24'
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("D11").Value = '15.86'
$ws.Range("E11").Value = '  +2.52%  '
$ws.Range("D12").Value = '0.06990'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = '0.6117'
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("D14").Value = '4.496'
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '0.9980'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '26.508.61'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '0.9975'
$ws.Range("D19").Value = '0.000007225'
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").Value = '11.40'
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("D21").Value = '1.947.48'
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = '4.472'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '8.558'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '5.101'
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("D25").Value = '138.18'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = '15.31'
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("D27").Value = '1.770'
$ws.Range("E27").Value = '  +2.60%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '1.380'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '106.30'
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = '3.926'
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("D31").Value = '0.07981'
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").Value = '3.682'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '0.04474'
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.610'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.001'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6235'
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9348'
$ws.Range("E37").Value = '  +2.92%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.042'
$ws.Range("E38").Value = '  +2.98%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.420'
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '1.000'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01515'
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.608'
$ws.Range("E42").Value = '  +3.66%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '99.40'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3857'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.905'
$ws.Range("E45").Value = '  +2.97%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1159'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05384'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '7.838'
$ws.Range("E48").Value = '  +1.99%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.25'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '51.76'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.233'
$ws.Range("E51").Value = '  -1.74%  '
